$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that must stay plain TEXT cells
# (matching the source export, which always writes Price as inlineStr).
# Assigning a bare numeric-looking string via .Value auto-converts the cell
# to a Number (COM type inference, just like typing it into Excel), so we
# briefly force Text format, assign the literal string, then restore the
# default "Normal" style so no stray style index is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.442.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.245.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.50%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.630"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.86"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.623"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0951"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.89%  "

$ws.Range("E13").Value = "  -0.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.861"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.250.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "42.341.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.31%  "

$ws.Range("E18").Value = "  +4.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +64.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.41%  "

$ws.Range("E24").Value = "  +3.55%  "

$ws.Range("E25").Value = "  +0.05%  "

$ws.Range("E26").Value = "  -1.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +23.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0818"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.119"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.57%  "

$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.126"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.44%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0315"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "13.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "64.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.202"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.72%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.102"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.71%  "

$ws.Range("E51").Value = "  +0.99%  "
